# Update "Enfermeria" sheet data: advance reporting month from
# 2023-08-01 (serial 45139) to 2023-09-01 (serial 45170) and refresh the
# associated Procedimiento / Cantidad figures for each Sede.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45170

# Row 2 - Bulevar / INYECTOLOGÍA
$ws.Cells.Item(2,2).Value = $newDate
$ws.Cells.Item(2,4).Value = 84

# Row 3 - Bulevar / TOMA DE EKG
$ws.Cells.Item(3,2).Value = $newDate
$ws.Cells.Item(3,4).Value = 131

# Row 4 - Bulevar / LAVADO DE OÍDOS
$ws.Cells.Item(4,2).Value = $newDate
$ws.Cells.Item(4,4).Value = 13

# Row 5 - Bulevar / RETIRO DE PUNTOS
$ws.Cells.Item(5,2).Value = $newDate
$ws.Cells.Item(5,4).Value = 2

# Row 6 - Bulevar / TOMA DE TENSIÓN ARTERIAL (was CURACIÓN)
$ws.Cells.Item(6,2).Value = $newDate
$ws.Cells.Item(6,3).Value = "TOMA DE TENSIÓN ARTERIAL"
$ws.Cells.Item(6,4).Value = 12

# Row 7 - Bulevar / GLUCOMETRÍA (was TOMA DE TENSIÓN ARTERIAL)
$ws.Cells.Item(7,2).Value = $newDate
$ws.Cells.Item(7,3).Value = "GLUCOMETRÍA"
$ws.Cells.Item(7,4).Value = 3

# Row 8 - San Martin / INYECTOLOGÍA (was Bulevar / GLUCOMETRÍA)
$ws.Cells.Item(8,1).Value = "San Martin"
$ws.Cells.Item(8,2).Value = $newDate
$ws.Cells.Item(8,3).Value = "INYECTOLOGÍA"
$ws.Cells.Item(8,4).Value = 89

# Row 9 - San Martin / LAVADO DE OÍDOS (was INYECTOLOGÍA)
$ws.Cells.Item(9,2).Value = $newDate
$ws.Cells.Item(9,3).Value = "LAVADO DE OÍDOS"
$ws.Cells.Item(9,4).Value = 13

# Row 10 - San Martin / TOMA DE EKG  (was LAVADO DE OÍDOS)
$ws.Cells.Item(10,2).Value = $newDate
$ws.Cells.Item(10,3).Value = "TOMA DE EKG "
$ws.Cells.Item(10,4).Value = 147

# Row 11 - San Martin / RETIRO DE PUNTOS (was TOMA DE EKG )
$ws.Cells.Item(11,2).Value = $newDate
$ws.Cells.Item(11,3).Value = "RETIRO DE PUNTOS"
$ws.Cells.Item(11,4).Value = 1

# Row 12 - San Martin / CURACIÓN (was RETIRO DE PUNTOS)
$ws.Cells.Item(12,2).Value = $newDate
$ws.Cells.Item(12,3).Value = "CURACIÓN"
$ws.Cells.Item(12,4).Value = 2

# Row 13 - San Martin / GLUCOMETRÍA (was CURACIÓN)
$ws.Cells.Item(13,2).Value = $newDate
$ws.Cells.Item(13,3).Value = "GLUCOMETRÍA"
$ws.Cells.Item(13,4).Value = 1

# Row 14 - Cartagena / INYECTOLOGÍA
$ws.Cells.Item(14,2).Value = $newDate
$ws.Cells.Item(14,4).Value = 90

# Row 15 - Cartagena / TOMA DE EKG
$ws.Cells.Item(15,2).Value = $newDate
$ws.Cells.Item(15,4).Value = 117

# Row 16 - Cartagena / TOMA DE TENSÓN ARTERIAL
$ws.Cells.Item(16,2).Value = $newDate
$ws.Cells.Item(16,4).Value = 110

# Row 17 - Cartagena / CURACIÓN (only date changes)
$ws.Cells.Item(17,2).Value = $newDate

# Row 18 - Cartagena / RETIRO DE PUNTOS
$ws.Cells.Item(18,2).Value = $newDate
$ws.Cells.Item(18,4).Value = 1

# Row 19 - Cartagena / GLUCOMETRÍA
$ws.Cells.Item(19,2).Value = $newDate
$ws.Cells.Item(19,4).Value = 4

# Update the active selection to D8, matching the author's final cursor
# position when they saved the workbook.
$ws.Range("D8").Select()
